$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- New label/value columns on Hoja2 (sheet2) ---
# Column G: necklace-length labels; Column I: size-code labels.
# Values are entered in the same order the shared-string table grows in the
# target workbook (G2..G6 first, then I6, then I2..I5).
$ws2.Range("G2").Value = "gargantilla: "
$ws2.Range("G3").Value = "princesa:"
$ws2.Range("G4").Value = "matine"
$ws2.Range("G5").Value = "opera:"
$ws2.Range("G6").Value = "cuerdalarga:"

$ws2.Range("I6").Value = "XG"
$ws2.Range("I2").Value = "AA"
$ws2.Range("I3").Value = "A"
$ws2.Range("I4").Value = "M"
$ws2.Range("I5").Value = "G"

# Match the formatting already used for the other label column (B): Segoe UI
# 8pt, left/center aligned with an indent, just in a darker grey tone.
$ws2.Range("B2").Copy()
$ws2.Range("G2:G6").PasteSpecial(-4122)
$ws2.Range("G2:G6").Font.Color = 855309

$ws2.Range("B2").Copy()
$ws2.Range("I2:I6").PasteSpecial(-4122)
$ws2.Range("I2:I6").Font.Color = 855309

$ws2.Columns.Item(7).ColumnWidth = 21.7265625

# --- Selection bookkeeping on Hoja3 (sheet3) ---
$ws3.Range("A1:A10").Select()

# --- Make Hoja2 the active sheet/tab, with the new size column selected ---
$ws2.Activate()
$ws2.Range("I2:I6").Select()
